$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DataRow($r, $d, $l, $m, $n, $o, $p, $q, $origin, $s, $t) {
    $ws.Cells.Item($r, 1).Value = 9
    $ws.Cells.Item($r, 2).Value = "Vega Central Mapocho de Santiago"
    $ws.Cells.Item($r, 3).Value = "Metropolitana"

    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 5).Value = 13
    $ws.Cells.Item($r, 6).Value = "Fruta"
    $ws.Cells.Item($r, 7).Value = 100104
    $ws.Cells.Item($r, 8).Value = "Frutos de pepita"
    $ws.Cells.Item($r, 9).Value = 100104004
    $ws.Cells.Item($r, 10).Value = "Níspero"
    $ws.Cells.Item($r, 11).Value = "Golden Nugget"

    $ws.Cells.Item($r, 12).Value = $l
    $ws.Cells.Item($r, 13).Value = $m
    $ws.Cells.Item($r, 14).Value = $n
    $ws.Cells.Item($r, 15).Value = $o
    $ws.Cells.Item($r, 16).Value = $p
    $ws.Cells.Item($r, 17).Value = $q
    $ws.Cells.Item($r, 18).Value = $origin
    $ws.Cells.Item($r, 19).Value = $s
    $ws.Cells.Item($r, 20).Value = $t
}

# Push the previously existing rows 2-5 down to rows 5-8 (unchanged data,
# preserved exactly as it was before this week's entries were added on top).
Set-DataRow 8 44915 "Primera"  200 5000 5000 5000 "`$/bandeja 5 kilos" "Provincia de Quillota" 1000 5
Set-DataRow 7 44915 "Especial" 150 6000 6000 6000 "`$/bandeja 5 kilos" "Provincia de Quillota" 1200 5
Set-DataRow 6 44911 "Segunda"  200 4000 4000 4000 "`$/bandeja 5 kilos" "Región de O'Higgins"   800  5
Set-DataRow 5 44911 "Primera"  220 5000 5000 5000 "`$/bandeja 5 kilos" "Región de O'Higgins"   1000 5

# Write this week's new data into rows 2-4.
Set-DataRow 2 45251 "Primera" 20 15000 15000 15000 "`$/bandeja 10 kilos" "Provincia de Quillota" 1500 10
Set-DataRow 3 45251 "Segunda" 40 12000 12000 12000 "`$/bandeja 10 kilos" "Provincia de Quillota" 1200 10
Set-DataRow 4 45251 "Tercera" 35 10000 10000 10000 "`$/bandeja 10 kilos" "Provincia de Quillota" 1000 10
